$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 content corrections:
#  B6 ("Is there jitter?") : "no" -> "yes"
#  C6 ("Is it stable?")    : "yes" -> "No"
$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = "No"

# Move the active selection to D6 (matches the saved selection state in the file)
$ws.Range("D6").Select()
